$d = $word.ActiveDocument

$replacements = @(
    @("757×5=3785", "296×9=2664"),
    @("223×4=892", "355×9=3195"),
    @("329×9=2961", "449×7=3143"),
    @("932×5=4660", "766×7=5362"),
    @("360×8=2880", "999×5=4995"),
    @("937×4=3748", "812×7=5684"),
    @("762×3=2286", "951×6=5706"),
    @("923×8=7384", "826×6=4956"),
    @("582×6=3492", "922×5=4610"),
    @("959×3=2877", "926×7=6482"),
    @("383×2=766", "847×4=3388"),
    @("518×8=4144", "423×7=2961"),
    @("106×9=954", "990×2=1980"),
    @("158×4=632", "485×2=970"),
    @("258×2=516", "695×9=6255"),
    @("569×2=1138", "703×7=4921"),
    @("673×7=4711", "230×2=460"),
    @("330×8=2640", "336×4=1344"),
    @("172×7=1204", "207×4=828"),
    @("818×4=3272", "884×8=7072"),
    @("986×6=5916", "520×4=2080"),
    @("453×6=2718", "478×4=1912"),
    @("969×5=4845", "369×9=3321"),
    @("440×8=3520", "935×9=8415"),
    @("149×6=894", "409×9=3681")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
